$wb = $excel.ActiveWorkbook

# --- Sheet: income_by_category ---
# Fix the Asian/Italian rows (previously swapped) and correct the income totals
# for every category to the right (recomputed) figures.
$wsCat = $wb.Worksheets.Item("income_by_category")

$wsCat.Range("A2").Value = "Italian"
$wsCat.Range("B2").Value = 2948
$wsCat.Range("C2").Value = "'" + '$49,462.70'

$wsCat.Range("A3").Value = "Asian"
$wsCat.Range("B3").Value = 3470
$wsCat.Range("C3").Value = "'" + '$46,720.65'

$wsCat.Range("C4").Value = "'" + '$34,796.80'

$wsCat.Range("C5").Value = "'" + '$28,237.75'

# --- Sheet: charges&passive ---
# Rename the mislabeled "chargers" header to "charges".
$wsCp = $wb.Worksheets.Item("charges&passive")
$wsCp.Range("G1").Value = "charges"

# --- Sheet: total_income ---
# Populate the totals: overall income, charges and passive income,
# and rename the "Total Income" header to "Total cache".
$wsTi = $wb.Worksheets.Item("total_income")

$wsTi.Range("B1").Value = "Total cache"

$wsTi.Range("B2").Value = 159217.9

$wsTi.Range("A3").Value = "charges"
$wsTi.Range("B3").Value = 71648.09

$wsTi.Range("A4").Value = "passive"
$wsTi.Range("B4").Value = 87569.89999999999
